$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Update fluid_mass (B11) which feeds Q_cool (B12 = fluid_mass*c_water*dt_cool)
$ws.Range("B11").Value = 0.58874952899999999

$wb.Save()
